# copy (not move) loc data from tourney sheets. in prep for code change.
#
# The "Tournament" sheet's `tournament` table gains two new key/value rows
# up top (competition-key, host-key) and eight new key/value rows at the
# bottom (venue-key.1 .. venue-key.8), copying identifying keys alongside
# the existing localized data rather than replacing it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tournament")
$lo = $ws.ListObjects.Item(1)

# Make room for the two new rows at the top (shifts existing rows 2-12
# down to rows 4-14, carrying their values/formatting along).
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "competition-key"
$ws.Range("B2").Value = "mens-world-cup"

$ws.Range("A3").Value = "host-key"
$ws.Range("B3").Value = "qatar"

# Append the eight venue-key rows after the (now shifted) last data row.
$venueKeys = @(
  @("venue-key.1", "qa-al-bayt"),
  @("venue-key.2", "qa-khalifa"),
  @("venue-key.3", "qa-al-thumama"),
  @("venue-key.4", "qa-ahmad-bin-ali"),
  @("venue-key.5", "qa-lusail"),
  @("venue-key.6", "qa-974"),
  @("venue-key.7", "qa-education-city"),
  @("venue-key.8", "qa-al-janoub")
)

$r = 15
foreach ($pair in $venueKeys) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $r = $r + 1
}

# Grow the table definition (and its autofilter) to cover the new rows.
$lo.Resize($ws.Range("A1:I22"))

# Match the editor's final selection state (the two newly-inserted rows).
$ws.Range("A2:XFD3").Select()
